$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rows for each team group (HẠNG n ĐỘI: CODE), followed by the two
# athlete rows belonging to that team.
$data = @(
    @{ Header = "HẠNG 1 ĐỘI: DHA"; HeaderRow = 3; Rows = @(
        @{ Row = 4;  A = 17; B = "Trần Vũ An Nhiên";     C = "DHA"; D = 2;  E = 5;   F = 6;  G = 9.5 },
        @{ Row = 5;  A = 12; B = "Phan Văn Gia Vũ";      C = "DHA"; D = 4;  E = 4.5 }
    ) },
    @{ Header = "HẠNG 2 ĐỘI: GLI"; HeaderRow = 6; Rows = @(
        @{ Row = 7;  A = 18; B = "Tạ Nguyễn Thiện Nhân"; C = "GLI"; D = 1;  E = 10;  F = 7;  G = 14 },
        @{ Row = 8;  A = 11; B = "Nguyễn Đức Duy";        C = "GLI"; D = 6;  E = 4 }
    ) },
    @{ Header = "HẠNG 3 ĐỘI: DKR"; HeaderRow = 9; Rows = @(
        @{ Row = 10; A = 15; B = "Trần Phúc Minh";        C = "DKR"; D = 3;  E = 4.5; F = 13; G = 8 },
        @{ Row = 11; A = 14; B = "Trần Bảo Nam";          C = "DKR"; D = 10; E = 3.5 }
    ) },
    @{ Header = "HẠNG 4 ĐỘI: HLA"; HeaderRow = 12; Rows = @(
        @{ Row = 13; A = 16; B = "Trần Tiến Hà";          C = "HLA"; D = 7;  E = 4;   F = 15; G = 8 },
        @{ Row = 14; A = 9;  B = "Nguyễn Việt Anh";       C = "HLA"; D = 8;  E = 4 }
    ) },
    @{ Header = "HẠNG 5 ĐỘI: CLO"; HeaderRow = 15; Rows = @(
        @{ Row = 16; A = 7;  B = "Nguyễn Thị Kim Anh";    C = "CLO"; D = 5;  E = 4.5; F = 16; G = 8 },
        @{ Row = 17; A = 1;  B = "Huỳnh Phúc Lâm";        C = "CLO"; D = 11; E = 3.5 }
    ) },
    @{ Header = "HẠNG 6 ĐỘI: HHO"; HeaderRow = 18; Rows = @(
        @{ Row = 19; A = 3;  B = "Hồ Nguyễn Vân Chi";     C = "HHO"; D = 12; E = 3.5; F = 28; G = 6.5 },
        @{ Row = 20; A = 8;  B = "Nguyễn Tùng Lâm";       C = "HHO"; D = 16; E = 3 }
    ) },
    @{ Header = "HẠNG 7 ĐỘI: TPH"; HeaderRow = 21; Rows = @(
        @{ Row = 22; A = 13; B = "Phạm Ngọc Dũng";        C = "TPH"; D = 13; E = 3.5; F = 28; G = 6.5 },
        @{ Row = 23; A = 2;  B = "Hồ Hùng Anh";           C = "TPH"; D = 15; E = 3 }
    ) },
    @{ Header = "HẠNG 8 ĐỘI: VLI"; HeaderRow = 24; Rows = @(
        @{ Row = 25; A = 10; B = "Nguyễn Đức Cao";        C = "VLI"; D = 17; E = 2.5; F = 35; G = 5 },
        @{ Row = 26; A = 4;  B = "Lê Bảo Tín";            C = "VLI"; D = 18; E = 2.5 }
    ) }
)

foreach ($group in $data) {
    $ws.Cells.Item($group.HeaderRow, 1).Value = $group.Header

    foreach ($r in $group.Rows) {
        $ws.Cells.Item($r.Row, 1).Value = $r.A
        $ws.Cells.Item($r.Row, 2).Value = $r.B
        $ws.Cells.Item($r.Row, 3).Value = $r.C
        $ws.Cells.Item($r.Row, 4).Value = $r.D
        $ws.Cells.Item($r.Row, 5).Value = $r.E
        if ($r.ContainsKey("F")) {
            $ws.Cells.Item($r.Row, 6).Value = $r.F
            $ws.Cells.Item($r.Row, 7).Value = $r.G
        }
    }
}
